$wb = $excel.ActiveWorkbook

# --- Sheet "LP1912": append 16 new scrape rows (1090-1105) ---
$ws1 = $wb.Worksheets.Item("LP1912")

$lp1912Rows = @(
    @("16:37:58", "16:51", "16_SANTA ANA", 14, "LP1912", "31/12/2025"),
    @("16:37:58", "17:03", "23_HERNANDEZ", 26, "LP1912", "31/12/2025"),
    @("16:37:58", "17:04", "14_ABASTO", 27, "LP1912", "31/12/2025"),
    @("16:37:58", "17:07", "15_ABASTO", 30, "LP1912", "31/12/2025"),
    @("16:37:58", "17:13", "10_OLMOS", 36, "LP1912", "31/12/2025"),
    @("16:37:58", "17:23", "16_SANTA ANA", 46, "LP1912", "31/12/2025"),
    @("16:37:58", "17:24", "11_ETCHEVERRY", 47, "LP1912", "31/12/2025"),
    @("16:37:58", "17:27", "15_ABASTO", 50, "LP1912", "31/12/2025"),
    @("16:37:58", "17:28", "23_HERNANDEZ", 51, "LP1912", "31/12/2025"),
    @("16:37:58", "17:34", "10_OLMOS", 57, "LP1912", "31/12/2025"),
    @("16:37:58", "17:35", "16_P MOR-SANTA ANA", 58, "LP1912", "31/12/2025"),
    @("16:37:58", "17:38", "17X38_ROMERO", 61, "LP1912", "31/12/2025"),
    @("16:37:58", "17:47", "16_SANTA ANA", 70, "LP1912", "31/12/2025"),
    @("16:37:58", "17:54", "10_OLMOS", 77, "LP1912", "31/12/2025"),
    @("16:37:58", "18:00", "23_HERNANDEZ", 83, "LP1912", "31/12/2025"),
    @("16:37:58", "18:04", "14_ABASTO", 87, "LP1912", "31/12/2025"),
)

$startRow = 1090
for ($i = 0; $i -lt $lp1912Rows.Count; $i++) {
    $r = $startRow + $i
    $data = $lp1912Rows[$i]
    $ws1.Cells.Item($r, 2).Value = $data[0]
    $ws1.Cells.Item($r, 3).Value = $data[1]
    $ws1.Cells.Item($r, 4).Value = $data[2]
    $ws1.Cells.Item($r, 5).Value = $data[3]
    $ws1.Cells.Item($r, 6).Value = $data[4]
    $ws1.Cells.Item($r, 7).Value = $data[5]
    $ws1.Cells.Item($r, 1).Style = "Normal"
}

$ws1.Range("A2").Value = "Última actualización: 31/12/2025 16:38:09"
$ws1.Range("A3").Value = "Total filas: 1104"

# --- Sheet "LP1912-215": timestamp only ---
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 31/12/2025 16:38:09"

# --- Sheet "6203-6173": append 1 new scrape row (136) + timestamp/total ---
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(136, 2).Value = "31/12/2025"
$ws3.Cells.Item(136, 3).Value = "16:38:03"
$ws3.Cells.Item(136, 4).Value = "16:58"
$ws3.Cells.Item(136, 5).Value = "215C_LA PLATA"
$ws3.Cells.Item(136, 6).Value = 20
$ws3.Cells.Item(136, 7).Value = "L6203"
$ws3.Cells.Item(136, 1).Style = "Normal"

$ws3.Range("A2").Value = "Última actualización: 31/12/2025 16:38:09"
$ws3.Range("A3").Value = "Total filas: 135"

